# edit.ps1 - apply "neutral voter language" + section reorder changes
# to dheeraj_chand_marketing_long_modern_clean.docx

$d = $word.ActiveDocument

function Find-ParaIndex {
    param($doc, [string]$matchText)
    $idx = 1
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$matchText*") {
            return $idx
        }
        $idx++
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Professional summary: plain text swap (no bold run needed here)
#    "...affecting all Black and Asian-American voters, developed..."
#    -> "...affecting 50M voters, developed..."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2)

# ---------------------------------------------------------------------------
# 2) Siege Analytics bullet: split the run so "50M" is bold + colored,
#    replacing "all Black and Asian-American" (note: "voters" text is kept
#    as-is, unbolded, immediately after).
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("all Black and Asian-American voters, developed geospatial machine")
if ($r.Find.Found) {
    $target = $d.Range($r.Start, $r.Start + 28)   # length of "all Black and Asian-American"
    $target.Text = "50M"
    $target.Font.Bold = 1
    $target.Font.Color = 5258796   # RGB 2C3E50 stored as BGR OLE color
}

# ---------------------------------------------------------------------------
# 3) Key Projects impact statement:
#    "...affecting all Black and Asian-American voters, improved..."
#    -> "...affecting 50M voters nationwide, improved..."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, improved electoral prediction accuracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters nationwide, improved electoral prediction accuracy", 2)

# ---------------------------------------------------------------------------
# 4) Move the "Analytics Supervisor - GSD&M" block (heading + 4 paragraphs)
#    so it appears after the "Data Products Manager - Helm/Murmuration"
#    block (i.e. right before "Senior Analyst - Myers Research").
# ---------------------------------------------------------------------------
$startIdx = Find-ParaIndex $d "Analytics Supervisor - GSD&M"
$endIdx   = Find-ParaIndex $d "Advanced Statistical and ML techniques for segmentation"

$moveStart = $d.Paragraphs($startIdx).Range.Start
$moveEnd   = $d.Paragraphs($endIdx).Range.End
$moveRange = $d.Range($moveStart, $moveEnd)
$moveRange.Cut()

$destIdx = Find-ParaIndex $d "Senior Analyst - Myers Research"
$destPara = $d.Paragraphs($destIdx)
$pasteRange = $d.Range($destPara.Range.Start, $destPara.Range.Start)
$pasteRange.Paste()

# The paste does not always carry over the paragraph style of the first
# moved paragraph (the "Analytics Supervisor..." heading) - restore it.
$headingIdx = Find-ParaIndex $d "Analytics Supervisor - GSD&M"
$headingPara = $d.Paragraphs($headingIdx)
if ($headingPara.Style.NameLocal -ne "Heading 3") {
    $headingPara.Style = "Heading 3"
}
